$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("F7").Value = 1624.15
$ws.Range("G7").Value = 1653.15
$ws.Range("H7").Value = 1620.6
$ws.Range("I7").Value = 1633.85
$ws.Range("J7").Value = 1646.95

# Row 9
$ws.Range("G9").Value = 1638
$ws.Range("H9").Value = 1620.4
$ws.Range("I9").Value = 1628.85

# Row 10
$ws.Range("G10").Value = 1641.25
$ws.Range("H10").Value = 1628
$ws.Range("I10").Value = 1640.85

# Row 11
$ws.Range("G11").Value = 1653.15
$ws.Range("H11").Value = 1639.75
$ws.Range("I11").Value = 1644.45

# Row 12
$ws.Range("G12").Value = 1652.3
$ws.Range("H12").Value = 1644
$ws.Range("I12").Value = 1648.95

# Row 13
$ws.Range("G13").Value = 1651
$ws.Range("H13").Value = 1642.5
$ws.Range("I13").Value = 1647.65

# Row 14
$ws.Range("G14").Value = 1647.9
$ws.Range("H14").Value = 1630.75
$ws.Range("I14").Value = 1631.25

# Row 15
$ws.Range("G15").Value = 1635
$ws.Range("H15").Value = 1628.05
$ws.Range("I15").Value = 1633.05

# Row 16
$ws.Range("G16").Value = 1636.75
$ws.Range("H16").Value = 1631.45
$ws.Range("I16").Value = 1635.65

# Row 17
$ws.Range("G17").Value = 1641.85
$ws.Range("H17").Value = 1634.4
$ws.Range("I17").Value = 1641.5

# Row 18
$ws.Range("G18").Value = 1644
$ws.Range("H18").Value = 1634.55
$ws.Range("I18").Value = 1635.8

# Row 19
$ws.Range("G19").Value = 1636.55
$ws.Range("H19").Value = 1626.8
$ws.Range("I19").Value = 1630.95

# Row 20
$ws.Range("G20").Value = 1636.9
$ws.Range("H20").Value = 1628.4
$ws.Range("I20").Value = 1633.7

# Row 21
$ws.Range("G21").Value = 1642
$ws.Range("H21").Value = 1629.2
$ws.Range("I21").Value = 1639
